$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is being restructured from "one row per person" (id, name,
# phone_number columns) to "one column per person" (Person1..Person5),
# with row labels id/name/phone and brand-new phone numbers. The old
# rows 5 and 6 are no longer part of the (now smaller) table, so clear
# them out entirely.
$ws.Range("A5:F6").Clear()

# Header row (B1:F1): Person1..Person5.
# B1 already carries the bold/border/center style (s=1) from the old
# "id" header, so reuse its formatting for the two brand-new header
# cells E1/F1 before writing any of the header values.
$ws.Range("B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

$ws.Range("B1").Value = "Person1"
$ws.Range("C1").Value = "Person2"
$ws.Range("D1").Value = "Person3"
$ws.Range("E1").Value = "Person4"
$ws.Range("F1").Value = "Person5"

# Row-label column (A2:A4) - already styled (s=1) from the previous
# "id" index column, so just overwrite the values.
$ws.Range("A2").Value = "id"
$ws.Range("A3").Value = "name"
$ws.Range("A4").Value = "phone"

# Row 2: ids
$ws.Range("B2").Value = 12345
$ws.Range("C2").Value = 123456
$ws.Range("D2").Value = 234567
$ws.Range("E2").Value = 345678
$ws.Range("F2").Value = 456789

# Row 3: names
$ws.Range("B3").Value = "Maxim"
$ws.Range("C3").Value = "Jasy"
$ws.Range("D3").Value = "Alex"
$ws.Range("E3").Value = "Adam"
$ws.Range("F3").Value = "Yura"

# Row 4: phone numbers
$ws.Range("B4").Value = "+375(44)6040934"
$ws.Range("C4").Value = "+375(44)9500689"
$ws.Range("D4").Value = "+375(44)7701371"
$ws.Range("E4").Value = "+375(44)4560025"
$ws.Range("F4").Value = "+375(44)2833246"
